# feat: add 2022-Q1 data
#
# - Duplicates the "2021-Q4" sheet (same fund-holding table layout) to
#   create a new "2022-Q1" sheet, positioned right after "2021-Q4" and
#   right before "总计".
# - Fills the new sheet with the 2022-Q1 fund data.
# - Inserts a new top data row in "总计" for 2022-Q1 and renumbers the
#   index column.

function Set-TextValue($range, [string]$text) {
    # Force the value to be stored as text (not auto-coerced to a number)
    # the same way the source workbook stores these numeric-looking
    # strings, then drop the helper number-format style again so the
    # cell doesn't pick up a stray style index.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet by copying "2021-Q4" (rId5) — it already
#    has the exact fund-holding table structure/styles we need — and
#    drop it into the same place, then rename + refill with new data.
# ------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Copy($null, $q4Sheet)

$newSheet = $wb.Worksheets.Item($q4Sheet.Index + 1)
$newSheet.Name = "2022-Q1"

# Row 2: 090019 / 大成景恒混合A
Set-TextValue $newSheet.Range("D2") "2.31"
Set-TextValue $newSheet.Range("E2") "93.51"
Set-TextValue $newSheet.Range("F2") "1.92"
Set-TextValue $newSheet.Range("G2") "0.0444"
$newSheet.Range("H2").Value = 8

# Row 3: 006038 / 大成景恒混合C
Set-TextValue $newSheet.Range("D3") "0.92"
Set-TextValue $newSheet.Range("E3") "93.51"
Set-TextValue $newSheet.Range("F3") "1.92"
Set-TextValue $newSheet.Range("G3") "0.0177"
$newSheet.Range("H3").Value = 8

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q1 at
#    the top of the data (row 2), shifting the rest down, then
#    renumber the index column (A) sequentially.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Give A2 the same style as the other index cells (copy format from A3).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Application.CutCopyMode = $false

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.06

# Renumber the index column sequentially (0,1,2,3,4,5) down through the
# row that used to be the last one (now shifted down by one).
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
